# Update column F ("dSF") values on the active worksheet to reflect
# the repulled / recalculated data from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 1
    5  = 6
    6  = 1
    7  = -1
    8  = -4
    9  = -1
    10 = -1
    12 = 4
    13 = 5
    15 = 2
    16 = 3
    17 = 3
    18 = -3
    19 = 4
    20 = 2
    21 = 1
    22 = 10
    24 = 3
    25 = -5
    26 = -2
    27 = -4
    28 = 7
    29 = 5
    30 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
